# DB_EJERCICIOS.xlsx - "Add files via upload"
# Populates the previously-empty Hoja1 sheet with the exercise database
# table (nomb / desc / tipo / imagen) and formats it the way the author's
# workbook ended up looking after the upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "nomb"
$ws.Range("B1").Value = "desc"
$ws.Range("C1").Value = "tipo"
$ws.Range("D1").Value = "imagen"

# --- Data, written in the same order the author originally typed it ---
$ws.Range("B2").Value = "Espalda pegada, bajar pecho"
$ws.Range("C2").Value = "Barra"
$ws.Range("A3").Value = "Peso Muerto"
$ws.Range("B3").Value = "Barra en trapecios, bajar rompiendo paralelo, espalda recta."
$ws.Range("A4").Value = "Overhead"
$ws.Range("B4").Value = "Empuje vertical estricto, bloqueo de codos arriba."
$ws.Range("A6").Value = "Remo Dorsal"
$ws.Range("A5").Value = "Zancada"
$ws.Range("A2").Value = "Press de Banca"
$ws.Range("A7").Value = "Squats"
$ws.Range("D7").Value = "sentadillabarra.jpg"
$ws.Range("D5").Value = "zancadabarra.jpg"
$ws.Range("D3").Value = "pesomuertobarra.jpg"
$ws.Range("D2").Value = "pressbancabarra.jpg"
$ws.Range("D4").Value = "overheadbarra.jpg"
$ws.Range("D6").Value = "remodorsalbarra.jpg"
$ws.Range("A8").Value = "Curl de bicesp"
$ws.Range("D8").Value = "curldebicepsbarra.jpg"

$ws.Range("C3").Value = "Barra"
$ws.Range("C4").Value = "Barra"
$ws.Range("C5").Value = "Barra"
$ws.Range("C6").Value = "Barra"
$ws.Range("C7").Value = "Barra"
$ws.Range("C8").Value = "Barra"

# --- Formatting: descriptions on rows 3-4 use a small Arial font ------
# (set per-cell, not as a combined range, so both cells resolve to the
# same single style index)
$ws.Range("B3").Font.Name = "Arial"
$ws.Range("B3").Font.Color = 1973274
$ws.Range("B4").Font.Name = "Arial"
$ws.Range("B4").Font.Color = 1973274

# --- Column widths ------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 24
$ws.Columns.Item(2).ColumnWidth = 62.166666666666664
$ws.Columns.Item(3).ColumnWidth = 21
$ws.Columns.Item(4).ColumnWidth = 27.166666666666668

# --- Page setup ---------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Leave the selection where the author left it ------------------------
$ws.Range("D12").Select() | Out-Null
